$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sprint 1 Daily SCRUM Standup")
$ws2 = $wb.Worksheets.Item("Sprint 2 Daily SCRUM Standup")

# --- Sprint 2 sheet: fill in the week's standup data ---

# Week label for this sheet's date range
$ws2.Range("C1").Value = "Week: 02/09/2025-02/15/2025"

# "What did you do yesterday?" block (rows 3-6)
$ws2.Range("B3").Value = "Whitley"
$ws2.Range("D3").Value = "Met with Dr. Bowman to discuss Sprint 1"
$ws2.Range("D3").WrapText = $true

$ws2.Range("B4").Value = "Asia"
$ws2.Range("D4").Value = "Met with Dr. Bowman to discuss Sprint 1"

$ws2.Range("B5").Value = "Erin"
$ws2.Range("D5").Value = "Met with Dr. Bowman to discuss Sprint 1"

$ws2.Range("B6").Value = "Jordan"
$ws2.Range("D6").Value = "Met with Dr. Bowman to discuss Sprint 1"
$ws2.Range("D6").WrapText = $true
$ws2.Range("D6").Borders.LineStyle = -4142

# "What will you work on today?" block (rows 7-10)
$ws2.Range("B7").Value = "Whitley"
$ws2.Range("D7").Value = "Project Plan, Form Project Report, Executive Summary"
$ws2.Range("D7").WrapText = $true

$ws2.Range("B8").Value = "Asia"
$ws2.Range("D8").Value = "Type up Executive Summary"
$ws2.Range("D8").WrapText = $true

$ws2.Range("B9").Value = "Erin"
$ws2.Range("D9").Value = "Fix Burn Down Chart, Update Project Report"

$ws2.Range("B10").Value = "Jordan"
$ws2.Range("D10").Value = "Put contraints in project report"
$ws2.Range("D10").WrapText = $true

# "Do you have any obstacles?" block (rows 11-14)
$ws2.Range("B11").Value = "Whitley"
$ws2.Range("D11").Value = "Limited availability"

$ws2.Range("B12").Value = "Asia"
$ws2.Range("D12").Value = "Work schedule conflicts"

$ws2.Range("B13").Value = "Erin"
$ws2.Range("D13").Value = "N/A"

$ws2.Range("B14").Value = "Jordan"
$b14Bottom = $ws2.Range("B14").Borders.Item(9)
$b14Bottom.LineStyle = 1
$b14Bottom.Weight = 2
$b14Bottom.Color = 0
$ws2.Range("D14").Value = "N/A"

# --- Window / selection state: move the active tab and selection from
#     Sprint 1 to Sprint 2, matching where the author left off editing ---
$ws1.Activate()
$ws1.Range("E3").Select()

$ws2.Activate()
$ws2.Range("D8").Select()
